$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '63.609.54'
$ws.Range("E2").Value = '  -3.16%  '

# Row 3
$ws.Range("D3").Value = '3.458.58'
$ws.Range("E3").Value = '  -3.63%  '

# Row 4
$ws.Range("E4").Value = '  -0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.35'
$ws.Range("E5").Value = '  -3.82%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.25'
$ws.Range("E6").Value = '  -5.32%  '

# Row 7
$ws.Range("D7").Value = '3.456.89'
$ws.Range("E7").Value = '  -3.69%  '

# Row 8
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.486'
$ws.Range("E9").Value = '  -2.37%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("E10").Value = '  -2.81%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.05'
$ws.Range("E11").Value = '  -2.86%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.382'
$ws.Range("E12").Value = '  -2.49%  '

# Row 13
$ws.Range("D13").Value = '4.048.70'
$ws.Range("E13").Value = '  -3.51%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.38'
$ws.Range("E14").Value = '  -2.79%  '

# Row 15
$ws.Range("E15").Value = '  -0.26%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000176'
$ws.Range("E16").Value = '  -5.99%  '

# Row 17
$ws.Range("D17").Value = '3.456.70'
$ws.Range("E17").Value = '  -3.68%  '

# Row 18
$ws.Range("D18").Value = '63.707.64'
$ws.Range("E18").Value = '  -3.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.82'
$ws.Range("E19").Value = '  -2.14%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.17'
$ws.Range("E20").Value = '  -3.24%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.61'
$ws.Range("E21").Value = '  -4.89%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.05'
$ws.Range("E22").Value = '  -2.08%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.570'
$ws.Range("E23").Value = '  -3.62%  '

# Row 24
$ws.Range("D24").Value = '3.601.53'
$ws.Range("E24").Value = '  -3.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '72.83'
$ws.Range("E25").Value = '  -1.99%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.03%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000107'
$ws.Range("E27").Value = '  -10.18%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.28%  '

# Row 29
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.53'
$ws.Range("E29").Value = '  -7.81%  '

# Row 30
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.36'
$ws.Range("E30").Value = '  -9.45%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.22'
$ws.Range("E31").Value = '  -7.91%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("E32").Value = '  -6.28%  '

# Row 33
$ws.Range("D33").Value = '3.461.48'
$ws.Range("E33").Value = '  -3.73%  '

# Row 34
$ws.Range("E34").Value = '  +0.01%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.48'
$ws.Range("E35").Value = '  -4.28%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.142'
$ws.Range("E36").Value = '  -3.64%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.19'
$ws.Range("E37").Value = '  -4.06%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '169.26'
$ws.Range("E38").Value = '  +0.46%  '

# Row 39
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.89'
$ws.Range("E39").Value = '  -2.87%  '

# Row 40
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.55'
$ws.Range("E40").Value = '  -4.00%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0797'
$ws.Range("E41").Value = '  -5.03%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.805'
$ws.Range("E42").Value = '  -4.50%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.12%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '25.45'
$ws.Range("E44").Value = '  -4.78%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.44'
$ws.Range("E45").Value = '  -3.97%  '

# Row 46
$ws.Range("E46").Value = '  -7.79%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.30'
$ws.Range("E47").Value = '  -5.46%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.61'
$ws.Range("E48").Value = '  -5.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.81'
$ws.Range("E49").Value = '  -3.37%  '

# Row 50
$ws.Range("D50").Value = '2.401.69'
$ws.Range("E50").Value = '  -1.99%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0263'
$ws.Range("E51").Value = '  -3.48%  '
